$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 111943947
$ws.Range("B3").Value = 85836
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 510
$ws.Range("F3").Value = "Doftskinn"
$ws.Range("G3").Value = "Cystostereum murrayi"
$ws.Range("H3").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q3").Value = 600352
$ws.Range("R3").Value = 7221402
$ws.Range("AX3").Value = "Maja Östlund, Simon Mattsson"

# Row 4
$ws.Range("A4").Value = 111943816
$ws.Range("B4").Value = 90816
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 2059
$ws.Range("F4").Value = "Skrovlig taggsvamp"
$ws.Range("G4").Value = "Hydnellum scabrosum"
$ws.Range("H4").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q4").Value = 600428
$ws.Range("R4").Value = 7221623
$ws.Range("AX4").Value = "Simon Mattsson, Maja Östlund"

# Row 5
$ws.Range("A5").Value = 111943877
$ws.Range("B5").Value = 90800
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q5").Value = 600476
$ws.Range("R5").Value = 7221499
$ws.Range("AX5").Value = "Simon Mattsson, Maja Östlund"

# Row 6
$ws.Range("A6").Value = 111943944
$ws.Range("B6").Value = 90800
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4364
$ws.Range("F6").Value = "Dropptaggsvamp"
$ws.Range("G6").Value = "Hydnellum ferrugineum"
$ws.Range("H6").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q6").Value = 600345
$ws.Range("R6").Value = 7221526
$ws.Range("AX6").Value = "Maja Östlund, Simon Mattsson"

# Row 7
$ws.Range("A7").Value = 111943940
$ws.Range("B7").Value = 77636
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 600237
$ws.Range("R7").Value = 7221447
$ws.Range("AX7").Value = "Maja Östlund, Simon Mattsson"

# Row 8
$ws.Range("A8").Value = 111943803
$ws.Range("B8").Value = 89539
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = "Ullticka"
$ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q8").Value = 600424
$ws.Range("R8").Value = 7221684
$ws.Range("AX8").Value = "Simon Mattsson, Maja Östlund"

# Row 9
$ws.Range("A9").Value = 111943887
$ws.Range("B9").Value = 90794
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 4362
$ws.Range("F9").Value = "Blå taggsvamp"
$ws.Range("G9").Value = "Hydnellum caeruleum"
$ws.Range("H9").Value = "(Hornem.) P.Karst."
$ws.Range("Q9").Value = 600485
$ws.Range("R9").Value = 7221470
$ws.Range("AX9").Value = "Simon Mattsson, Maja Östlund"

# Row 10
$ws.Range("A10").Value = 111943879
$ws.Range("B10").Value = 90800
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 4364
$ws.Range("F10").Value = "Dropptaggsvamp"
$ws.Range("G10").Value = "Hydnellum ferrugineum"
$ws.Range("H10").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q10").Value = 600452
$ws.Range("R10").Value = 7221545
$ws.Range("AX10").Value = "Simon Mattsson, Maja Östlund"

# Row 11
$ws.Range("A11").Value = 111943881
$ws.Range("B11").Value = 90800
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 4364
$ws.Range("F11").Value = "Dropptaggsvamp"
$ws.Range("G11").Value = "Hydnellum ferrugineum"
$ws.Range("H11").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q11").Value = 600419
$ws.Range("R11").Value = 7221630
$ws.Range("AX11").Value = "Simon Mattsson, Maja Östlund"

# Row 12
$ws.Range("A12").Value = 111943814
$ws.Range("B12").Value = 90816
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 2059
$ws.Range("F12").Value = "Skrovlig taggsvamp"
$ws.Range("G12").Value = "Hydnellum scabrosum"
$ws.Range("H12").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q12").Value = 600437
$ws.Range("R12").Value = 7221630
$ws.Range("AX12").Value = "Simon Mattsson, Maja Östlund"

# Row 13
$ws.Range("A13").Value = 111943907
$ws.Range("B13").Value = 90816
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 2059
$ws.Range("F13").Value = "Skrovlig taggsvamp"
$ws.Range("G13").Value = "Hydnellum scabrosum"
$ws.Range("H13").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q13").Value = 600409
$ws.Range("R13").Value = 7221648
$ws.Range("AX13").Value = "Maja Östlund, Simon Mattsson"

# Row 14
$ws.Range("A14").Value = 111943883
$ws.Range("B14").Value = 90800
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = "Dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum ferrugineum"
$ws.Range("H14").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q14").Value = 600311
$ws.Range("R14").Value = 7221358
$ws.Range("AX14").Value = "Simon Mattsson, Maja Östlund"

# Row 15
$ws.Range("A15").Value = 111943815
$ws.Range("B15").Value = 90816
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 2059
$ws.Range("F15").Value = "Skrovlig taggsvamp"
$ws.Range("G15").Value = "Hydnellum scabrosum"
$ws.Range("H15").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q15").Value = 600430
$ws.Range("R15").Value = 7221629
$ws.Range("AX15").Value = "Simon Mattsson, Maja Östlund"

# Row 16
$ws.Range("A16").Value = 111943882
$ws.Range("B16").Value = 90800
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4364
$ws.Range("F16").Value = "Dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum ferrugineum"
$ws.Range("H16").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q16").Value = 600419
$ws.Range("R16").Value = 7221432
$ws.Range("AX16").Value = "Simon Mattsson, Maja Östlund"

# Row 17
$ws.Range("A17").Value = 111943880
$ws.Range("B17").Value = 90800
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 4364
$ws.Range("F17").Value = "Dropptaggsvamp"
$ws.Range("G17").Value = "Hydnellum ferrugineum"
$ws.Range("H17").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q17").Value = 600437
$ws.Range("R17").Value = 7221453
$ws.Range("AX17").Value = "Simon Mattsson, Maja Östlund"

# Row 18
$ws.Range("A18").Value = 111943841
$ws.Range("B18").Value = 77389
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 228912
$ws.Range("F18").Value = "Mörk kolflarnlav"
$ws.Range("G18").Value = "Carbonicola myrmecina"
$ws.Range("H18").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q18").Value = 600367
$ws.Range("R18").Value = 7221297
$ws.Range("AX18").Value = "Simon Mattsson, Maja Östlund"

